$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.336.14"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -3.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.563.92"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -4.18%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.31"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -5.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.66"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.559.70"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -4.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.612"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.04%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.668"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -6.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.145"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -10.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "52.71"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -7.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000259"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -10.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.77"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -7.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.133.64"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.565.38"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.19%  "
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.32"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -5.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.17"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -6.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "66.236.87"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.05"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -7.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "394.07"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.31"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -6.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.82"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.19"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.89"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -5.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.43"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.04"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.54"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -6.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.92"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -7.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.97"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -6.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.03"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.14"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "619.87"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.83%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "63.58"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.65%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.113"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -8.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "41.24"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -7.54%  "
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.394"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0760"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -8.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.131"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -6.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.991.36"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.82"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -7.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.51"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0408"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -8.18%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.12"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.130"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -6.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.50"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -6.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.12"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.71"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.05%  "
